$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.014929244053271
$ws.Cells.Item(2, 4).Value = 1.020775119883531
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.013247807617859
$ws.Cells.Item(2, 9).Value = 1.025851439972041
$ws.Cells.Item(2, 10).Value = 1.020157669769753
$ws.Cells.Item(2, 11).Value = 1.023614408938973
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.016109531397283
$ws.Cells.Item(2, 14).Value = 1.010807519146486
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.016205220458409
$ws.Cells.Item(3, 4).Value = 1.021692452815325
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.015173191787953
$ws.Cells.Item(3, 9).Value = 1.026105573556114
$ws.Cells.Item(3, 10).Value = 1.021066714046317
$ws.Cells.Item(3, 11).Value = 1.024337748708014
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.017836484460037
$ws.Cells.Item(3, 14).Value = 1.011110085313937
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.017028740503584
$ws.Cells.Item(4, 4).Value = 1.022283943223088
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.016416475128022
$ws.Cells.Item(4, 9).Value = 1.026267290172515
$ws.Cells.Item(4, 10).Value = 1.021652396820202
$ws.Cells.Item(4, 11).Value = 1.024803096572644
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.018950994704106
$ws.Cells.Item(4, 14).Value = 1.011304942843843
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.017374447098387
$ws.Cells.Item(5, 4).Value = 1.022532110353021
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.016938552299475
$ws.Cells.Item(5, 9).Value = 1.026334625025463
$ws.Cells.Item(5, 10).Value = 1.021898017334055
$ws.Cells.Item(5, 11).Value = 1.024998085784389
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.019418844274004
$ws.Cells.Item(5, 14).Value = 1.011386641646181
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.01743246356005
$ws.Cells.Item(6, 4).Value = 1.022573749762651
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.017026176678171
$ws.Cells.Item(6, 9).Value = 1.026345892722902
$ws.Cells.Item(6, 10).Value = 1.021939223030545
$ws.Cells.Item(6, 11).Value = 1.025030787729124
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.019497358193067
$ws.Cells.Item(6, 14).Value = 1.011400346422619
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.01703336181322
$ws.Cells.Item(7, 4).Value = 1.022287261186556
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.016423453483506
$ws.Cells.Item(7, 9).Value = 1.026268192459562
$ws.Cells.Item(7, 10).Value = 1.021655681166991
$ws.Cells.Item(7, 11).Value = 1.02480570455026
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.018957248828864
$ws.Cells.Item(7, 14).Value = 1.011306035366916
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.015360909732708
$ws.Cells.Item(8, 4).Value = 1.0210855709289
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.013899041726932
$ws.Cells.Item(8, 9).Value = 1.025937890715056
$ws.Cells.Item(8, 10).Value = 1.020465412398501
$ws.Cells.Item(8, 11).Value = 1.023859426557731
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.016693781419973
$ws.Cells.Item(8, 14).Value = 1.010909964925029
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.012397276805133
$ws.Cells.Item(9, 4).Value = 1.018951895106328
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.009430315145557
$ws.Cells.Item(9, 9).Value = 1.025334918724303
$ws.Cells.Item(9, 10).Value = 1.018348402346849
$ws.Cells.Item(9, 11).Value = 1.02217110768792
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.012682074527158
$ws.Cells.Item(9, 14).Value = 1.010204893458853
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.01040995502657
$ws.Cells.Item(10, 4).Value = 1.017518343460313
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.00643644980466
$ws.Cells.Item(10, 9).Value = 1.024918759703272
$ws.Cells.Item(10, 10).Value = 1.016923555091584
$ws.Cells.Item(10, 11).Value = 1.02103130101419
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.009991111079376
$ws.Cells.Item(10, 14).Value = 1.009729934948374
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.009546581542897
$ws.Cells.Item(11, 4).Value = 1.016894910330066
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.005136358141494
$ws.Cells.Item(11, 9).Value = 1.024735172282214
$ws.Cells.Item(11, 10).Value = 1.016303302858407
$ws.Cells.Item(11, 11).Value = 1.02053431530909
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.008821780499811
$ws.Cells.Item(11, 14).Value = 1.009523083597871
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.009225449441732
$ws.Cells.Item(12, 4).Value = 1.016662929482969
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.004652866611126
$ws.Cells.Item(12, 9).Value = 1.024666468649492
$ws.Cells.Item(12, 10).Value = 1.016072413681836
$ws.Cells.Item(12, 11).Value = 1.020349190758286
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.008386801712228
$ws.Cells.Item(12, 14).Value = 1.009446068696294
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.00929435326985
$ws.Cells.Item(13, 4).Value = 1.016712708797869
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.00475660372884
$ws.Cells.Item(13, 9).Value = 1.024681228974122
$ws.Cells.Item(13, 10).Value = 1.016121962919925
$ws.Cells.Item(13, 11).Value = 1.020388924299624
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.008480135263294
$ws.Cells.Item(13, 14).Value = 1.009462596883903
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.009520045622657
$ws.Cells.Item(14, 4).Value = 1.016875743114511
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.005096404517915
$ws.Cells.Item(14, 9).Value = 1.024729503656933
$ws.Cells.Item(14, 10).Value = 1.016284227722378
$ws.Cells.Item(14, 11).Value = 1.020519023531404
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.008785838121243
$ws.Cells.Item(14, 14).Value = 1.009516721228939
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.009659043998797
$ws.Cells.Item(15, 4).Value = 1.016976139454397
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.005305689609723
$ws.Cells.Item(15, 9).Value = 1.02475917949407
$ws.Cells.Item(15, 10).Value = 1.016384137987579
$ws.Cells.Item(15, 11).Value = 1.020599112655256
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.008974106759552
$ws.Cells.Item(15, 14).Value = 1.00955004496228
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.010467193505786
$ws.Cells.Item(16, 4).Value = 1.017559661440585
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.006522652399925
$ws.Cells.Item(16, 9).Value = 1.024930872239245
$ws.Cells.Item(16, 10).Value = 1.016964649447174
$ws.Cells.Item(16, 11).Value = 1.021064211391677
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.010068627285868
$ws.Cells.Item(16, 14).Value = 1.009743637712563
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.010973355292245
$ws.Cells.Item(17, 4).Value = 1.017924964109755
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.007285009176453
$ws.Cells.Item(17, 9).Value = 1.025037662070649
$ws.Cells.Item(17, 10).Value = 1.017327905174598
$ws.Cells.Item(17, 11).Value = 1.021355030664646
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.010754074240755
$ws.Cells.Item(17, 14).Value = 1.009864752916302
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.011268316495672
$ws.Cells.Item(18, 4).Value = 1.01813777921336
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.007729320594471
$ws.Cells.Item(18, 9).Value = 1.025099624042284
$ws.Cells.Item(18, 10).Value = 1.01753946968356
$ws.Cells.Item(18, 11).Value = 1.021524328864084
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.011153486980704
$ws.Cells.Item(18, 14).Value = 1.009935282579985
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.011368844357191
$ws.Cells.Item(19, 4).Value = 1.018210299753769
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.007880759102468
$ws.Cells.Item(19, 9).Value = 1.025120696117237
$ws.Cells.Item(19, 10).Value = 1.017611554240094
$ws.Cells.Item(19, 11).Value = 1.021581999046995
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.011289609613677
$ws.Cells.Item(19, 14).Value = 1.009959311967201
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.010919077369308
$ws.Cells.Item(20, 4).Value = 1.017885797496505
$ws.Cells.Item(20, 5).Value = 0.9894336180360677
$ws.Cells.Item(20, 6).Value = 1.007203252723311
$ws.Cells.Item(20, 9).Value = 1.025026238343461
$ws.Cells.Item(20, 10).Value = 1.017288964049144
$ws.Cells.Item(20, 11).Value = 1.021323862860701
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.010680573449354
$ws.Cells.Item(20, 14).Value = 1.009851770289901
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.009453596975634
$ws.Cells.Item(21, 4).Value = 1.01682774493108
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.004996357812175
$ws.Cells.Item(21, 9).Value = 1.024715302091788
$ws.Cells.Item(21, 10).Value = 1.016236458634189
$ws.Cells.Item(21, 11).Value = 1.020480726997275
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.008695833984539
$ws.Cells.Item(21, 14).Value = 1.009500787974899
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.008529660085026
$ws.Cells.Item(22, 4).Value = 1.016160130477646
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.003605432194656
$ws.Cells.Item(22, 9).Value = 1.024516846044956
$ws.Cells.Item(22, 10).Value = 1.015571811268992
$ws.Cells.Item(22, 11).Value = 1.019947591836609
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.007444253201973
$ws.Cells.Item(22, 14).Value = 1.009279062594915
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.009019699325878
$ws.Cells.Item(23, 4).Value = 1.016514272328448
$ws.Cells.Item(23, 5).Value = 0.9879432794636464
$ws.Cells.Item(23, 6).Value = 1.004343113534538
$ws.Cells.Item(23, 9).Value = 1.024622332456818
$ws.Cells.Item(23, 10).Value = 1.015924430059674
$ws.Cells.Item(23, 11).Value = 1.020230504970419
$ws.Cells.Item(23, 12).Value = 0.9917760702887611
$ws.Cells.Item(23, 13).Value = 1.008108096122234
$ws.Cells.Item(23, 14).Value = 1.009396703555828
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.010943604073403
$ws.Cells.Item(24, 4).Value = 1.017903496005337
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.007240196050603
$ws.Cells.Item(24, 9).Value = 1.02503140124341
$ws.Cells.Item(24, 10).Value = 1.017306560844591
$ws.Cells.Item(24, 11).Value = 1.021337947274362
$ws.Cells.Item(24, 12).Value = 0.9929938892766441
$ws.Cells.Item(24, 13).Value = 1.010713786517576
$ws.Cells.Item(24, 14).Value = 1.009857636934416
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.013165453533465
$ws.Cells.Item(25, 4).Value = 1.019505439748327
$ws.Cells.Item(25, 5).Value = 0.9912096547607051
$ws.Cells.Item(25, 6).Value = 1.010588107250644
$ws.Cells.Item(25, 9).Value = 1.025493292879149
$ws.Cells.Item(25, 10).Value = 1.018898056729389
$ws.Cells.Item(25, 11).Value = 1.022610074431632
$ws.Cells.Item(25, 12).Value = 0.9944092447426416
$ws.Cells.Item(25, 13).Value = 1.013722037066177
$ws.Cells.Item(25, 14).Value = 1.010388028962082
